# Regenerate save_data column G ("K") values for rows 2-15.
# This mirrors the source pipeline change that computes K (strikeouts-derived
# stat) differently than the old "Strike#" values, so only the literal
# numeric values in column G need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 2
    14 = 2
    15 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
